$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New catalog entries appended below the existing 42 comarques rows
# (row 46: Codi 98 / "No consta", row 47: Codi 99 / "Altres/Diversos").
# Format the cells as Text first so the numeric-looking codes are stored
# as text (matching the rest of column A), then write column A before
# column B so new shared-string entries are created in "98", "99",
# "No consta", "Altres/Diversos" order.
$ws.Range("A46:B47").NumberFormat = "@"

$ws.Range("A46").Value = "98"
$ws.Range("A47").Value = "99"
$ws.Range("B46").Value = "No consta"
$ws.Range("B47").Value = "Altres/Diversos"

# Grow the table ("Tabla1") so the new rows are included in it.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:B47"))
